# aktualizace vzoru, aby datumy byly v jednom mesici
# Shift the sample dates in the "Spoty" sheet forward by one month (30 days)
# so that every booking in the template falls inside a single calendar month.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spoty")

# column F (and G where present) hold serial date values; bump each one
# from its old value to the new one exactly as recorded in the template update
$ws.Range("F3").Value  = 45127
$ws.Range("F4").Value  = 45128
$ws.Range("F5").Value  = 45129
$ws.Range("F6").Value  = 45127
$ws.Range("F7").Value  = 45128
$ws.Range("F8").Value  = 45130
$ws.Range("F9").Value  = 45126
$ws.Range("F10").Value = 45127
$ws.Range("F11").Value = 45129

$ws.Range("F12").Value = 45119
$ws.Range("G12").Value = 45127

$ws.Range("F13").Value = 45110
$ws.Range("G13").Value = 45110

$ws.Range("F14").Value = 45127
$ws.Range("G14").Value = 45129
